# GreenLeaf_EVM.xlsx - Status Report 3 + EVM Update
# Updates the milestone dates on the "tabella" sheet header row and the
# Actual Cost / %Progress figures for the "30/12/22" milestone column.
# All dependent formulas (CV, CPI, ETC, EAC, VAC, Average Index, Stato)
# and the chart caches that read from this sheet recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tabella")

# Rename the milestone date headers (row 3, columns E:H).
# These also drive the Tabella1 table column names and the chart category
# (c:cat) caches, since both reference tabella!$B$3:$H$3 / the table header.
$ws.Range("E3").Value = "30/12/22"
$ws.Range("F3").Value = "23/01/23"
$ws.Range("G3").Value = "25/01/23"
$ws.Range("H3").Value = "26/01/23"

# Actual Cost (AC) for the "30/12/22" milestone was 0 (not yet started) and
# is now 43.75.
$ws.Range("E6").Value = 43.75

# %Progress for the same milestone moves from 0 to 1 (complete).
$ws.Range("E8").Value = 1
